# Feria Lagunitas de Puerto Montt - Betarraga: add a new weekly price record.
# A new row is inserted at row 483, pushing the existing rows 483:576 down
# to 484:577 (dimension grows from A1:R576 to A1:R577).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 483 (shifts rows 483:576 -> 484:577)
$ws.Rows(483).Insert()

# Populate the newly inserted row 483 with the new record
$ws.Range("A483").Value = 4
$ws.Range("B483").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C483").Value = "Los Lagos"
$ws.Range("D483").Value = 45275
$ws.Range("E483").Value = 10
$ws.Range("F483").Value = 100114014
$ws.Range("G483").Value = "Betarraga"
$ws.Range("H483").Value = "Sin especificar"
$ws.Range("I483").Value = "Primera"
$ws.Range("J483").Value = 1200
$ws.Range("K483").Value = 1100
$ws.Range("L483").Value = 1100
$ws.Range("M483").Value = 1100
$ws.Range("N483").Value = "$/paquete 5 unidades"
$ws.Range("O483").Value = "Región Metropolitana"
$ws.Range("P483").Value = 220
$ws.Range("Q483").Value = 5
$ws.Range("R483").Value = "Hortaliza"
